# Apply the "add solved solution and result table" edit to the TSP result
# header workbook: update the COPT remark note, refresh the result values
# for the re-run solves, resize the wrapped-text rows to a uniform height,
# and restore the window/selection state left behind by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# ---------------------------------------------------------------------
# 1. Remark (column J) text for the COPT rows (J3, J5, J7 all share the
#    same note) - note that the subtour-shortest-cycle function now also
#    breaks symmetry, and that a nearest-neighbour warm start is used
#    before handing the problem to Gurobi.
# ---------------------------------------------------------------------
$coptNote = "1. Callback: 只有找到合法子环，才添加lazy constraint把该子环排除；放弃在硬约束中排除所有子环`n2. 求最短子环的函数不会找到对称子环，以破解对称性`n3.在Gurobi求解之前使用最近邻算法生成初始可行路径，并设置为初始解`n4. 参数调优：`nmodel.setParam(COPT.Param.HeurLevel, 3)`n        model.setParam(COPT.Param.Threads, 8)`n        model.setParam(COPT.Param.Presolve, 2)"

$ws.Range("J3").Value2 = $coptNote
$ws.Range("J5").Value2 = $coptNote
$ws.Range("J7").Value2 = $coptNote

# ---------------------------------------------------------------------
# 2. Refresh the solved results (columns F-I) for the re-run instances.
# ---------------------------------------------------------------------
$ws.Range("F2").Value2 = 48832

$ws.Range("F3").Value2 = 47721.110999999997
$ws.Range("G3").Value2 = 59714
$ws.Range("H3").Value2 = 20.083881315749998
$ws.Range("I3").Value2 = 600.17999999999995

$ws.Range("F4").Value2 = 6770
$ws.Range("H4").Value2 = 0.044293158300000003
$ws.Range("I4").Value2 = 600.04200000000003

$ws.Range("F5").Value2 = 6761
$ws.Range("G5").Value2 = 8536
$ws.Range("H5").Value2 = 20.783193700000002
$ws.Range("I5").Value2 = 600.33000000000004

$ws.Range("I6").Value2 = 600.09

$ws.Range("F7").Value2 = 211007.25
$ws.Range("G7").Value2 = 305404
$ws.Range("H7").Value2 = 30.9088125
$ws.Range("I7").Value2 = 600.9

# ---------------------------------------------------------------------
# 3. All wrapped-text rows (2-7) now share one uniform custom height.
# ---------------------------------------------------------------------
2..7 | ForEach-Object {
    $ws.Rows.Item($_).RowHeight = 250.95
}

# ---------------------------------------------------------------------
# 4. Restore the view state: scroll position, zoom level and the last
#    active selection (topLeftCell A7, zoom 85%, selection F7).
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 85
$ws.Range("F7").Select() | Out-Null
